$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for column D (same text as B1: "CBi")
$ws.Range("D1").Value = "CBi"

# New values for D2:D93 (product category "CBi" re-run, slightly different precision)
$dValues = @(
    47.189413612031601,58.738083360244502,44.281694680581701,112.046124424308,39.333205569688097,
    94.415952422157901,81.0860457226687,53.282041561418701,135.68908554839101,51.479927151914403,
    53.395830584004301,50.748768676098699,386.65901138684598,164.354662761198,312.687602219349,
    95.374326539297201,41.067977955707804,102.64437028577299,45.926158487472101,46.564118020385102,
    35.137031853292001,43.983733953707201,218.52180976469799,138.20530932192099,28.6674775987745,
    106.234941636199,368.74045422218597,64.518981971608795,797.75015076397005,100.877668986823,
    109.421496702224,144.33575947774699,48.874210485486898,143.78196519989399,51.388989757533402,
    86.706898898595,89.341393943718501,102.404210675996,116.937246126544,39.764176175339998,
    182.625351657198,102.91848913702,104.668229073487,52.148399168236402,131.16814949537701,
    228.712823602306,89.568262540305994,134.729214256266,160.95349421880101,89.840918428082603,
    127.446315271572,279.60754606795598,457.23736602569699,0,367.44215918850301,
    1195.7385072396,11.060972679280701,601.69460500270702,439.54417329655399,883.80411528027605,
    57.214102174169398,279.77315035835198,36.015473032969197,12.3455593208305,15.1126225541346,
    19.988445783934601,39.269256088968397,46.366184165177103,44.3114391283356,150.368250150981,
    28327.550188678499,110.370512301864,41.2872142595266,20.406128906071501,13.186696677522001,
    19.6321071039605,46.826601475067498,70.411624441789598,117.48087767625201,31.160322645503701,
    277.86045513357698,133.28489008693299,105.841415679091,77.0756578571971,0,
    69.210845057260101,64.252323390676906,624.05146076423205,173.10718845322401,23.869309194377301,
    3.2020811100272701,0
)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Update the view state to match the edited workbook
$ws.Application.ActiveWindow.ScrollRow = 74
$ws.Range("B85:D86").Select()
$excel.ActiveWindow.RangeSelection.Item(2).Activate()
